# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect newly scraped counts (gh-pages output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll     = $wb.Worksheets.Item("全部类型")

# Row -> new F value for 展览 sheet
$exhibitUpdates = @{
    2  = 3074
    3  = 486
    5  = 49
    9  = 1057
    10 = 14906
    11 = 183
    12 = 144
    13 = 503
    14 = 5929
    22 = 199
    23 = 818
    24 = 2957
    26 = 10751
    28 = 82
}

foreach ($row in $exhibitUpdates.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitUpdates[$row]
}

# Row -> new F value for 全部类型 sheet
$allTypesUpdates = @{
    3  = 3074
    4  = 486
    6  = 49
    10 = 1057
    11 = 14906
    12 = 183
    13 = 144
    14 = 503
    15 = 5929
    23 = 199
    24 = 818
    25 = 2957
    28 = 10751
    30 = 82
}

foreach ($row in $allTypesUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allTypesUpdates[$row]
}
